$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row: "_old" -> "_FV2410", "_new" -> "_FV2504"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -like "*_old") {
        $cell.Value = ($v -replace "_old$", "_FV2410")
    } elseif ($v -like "*_new") {
        $cell.Value = ($v -replace "_new$", "_FV2504")
    }
}

# 2) Turn the used range into an Excel Table (adds xl/tables/table1.xml,
#    the worksheet rels and the <tableParts> element).
$rng = $ws.Range("A1:U82")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"

# 3) Freeze the header row (pane split after row 1).
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
